$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1875.5
$ws.Range("J32").Value = 1875.5
$ws.Range("L32").Value = 1875.5
$ws.Range("N32").Value = -2527.5
$ws.Range("H43").Value = 9165.846
$ws.Range("I43").Value = 2258.3333
$ws.Range("J43").Value = 15086.571
$ws.Range("K43").Value = 2258.3333
$ws.Range("L43").Value = 15086.571
$ws.Range("M43").Value = -2189.3333
$ws.Range("N43").Value = -15224.571
$ws.Range("H116").Value = 2476.8235
$ws.Range("I116").Value = 2250
$ws.Range("J116").Value = 3021.2
$ws.Range("K116").Value = 2250
$ws.Range("L116").Value = 3021.2
$ws.Range("M116").Value = 1192
$ws.Range("N116").Value = -9905.200000000001
$ws.Range("H132").Value = 5437748.5
$ws.Range("I132").Value = 3143.838
$ws.Range("J132").Value = 27780012
$ws.Range("K132").Value = 9431.514000000001
$ws.Range("L132").Value = 83340036
$ws.Range("M132").Value = -6901.514000000001
$ws.Range("N132").Value = -83345096
$ws.Range("H137").Value = 29019
$ws.Range("I137").Value = 1564
$ws.Range("J137").Value = 70201.5
$ws.Range("K137").Value = 4692
$ws.Range("L137").Value = 210604.5
$ws.Range("M137").Value = -2142
$ws.Range("N137").Value = -215704.5
$ws.Range("H141").Value = 1469.7059
$ws.Range("I141").Value = 1311.5625
$ws.Range("J141").Value = 4000
$ws.Range("K141").Value = 3934.6875
$ws.Range("L141").Value = 12000
$ws.Range("M141").Value = 1245.3125
$ws.Range("N141").Value = -22360

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2445.848
$ws.Range("I61").Value = 2268.1162
$ws.Range("J61").Value = 4993.3335
$ws.Range("K61").Value = 2268.1162
$ws.Range("L61").Value = 4993.3335
$ws.Range("M61").Value = -2056.1162
$ws.Range("N61").Value = -5417.3335
$ws.Range("H102").Value = 1500
$ws.Range("I102").Value = 1500
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 1500
$ws.Range("L102").Value = 0
$ws.Range("M102").ClearContents()
$ws.Range("N102").Value = 122
$ws.Range("H128").Value = 16679
$ws.Range("J128").Value = 16679
$ws.Range("L128").Value = 16679
$ws.Range("N128").Value = -26639
$ws.Range("H132").Value = 1731.1698
$ws.Range("I132").Value = 1288.475
$ws.Range("J132").Value = 3093.3076
$ws.Range("K132").Value = 3865.425
$ws.Range("L132").Value = 9279.9228
$ws.Range("M132").Value = -1335.425
$ws.Range("N132").Value = -14339.9228
$ws.Range("H136").Value = 2445.848
$ws.Range("I136").Value = 2268.1162
$ws.Range("J136").Value = 4993.3335
$ws.Range("K136").Value = 6804.348599999999
$ws.Range("L136").Value = 14980.0005
$ws.Range("M136").Value = -4254.348599999999
$ws.Range("N136").Value = -20080.0005

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2686.4546
$ws.Range("I86").Value = 2446.3333
$ws.Range("J86").Value = 3201
$ws.Range("K86").Value = 2446.3333
$ws.Range("L86").Value = 3201
$ws.Range("M86").Value = -1323.3333
$ws.Range("N86").Value = -5447
$ws.Range("H89").Value = 2686.4546
$ws.Range("I89").Value = 2446.3333
$ws.Range("J89").Value = 3201
$ws.Range("K89").Value = 12231.6665
$ws.Range("L89").Value = 16005
$ws.Range("M89").Value = -6615.666499999999
$ws.Range("N89").Value = -27237
$ws.Range("H118").Value = 7894.4443
$ws.Range("J118").Value = 7894.4443
$ws.Range("L118").Value = 7894.4443
$ws.Range("N118").Value = -11208.4443
$ws.Range("H134").Value = 2783.1025
$ws.Range("I134").Value = 1831.8823
$ws.Range("J134").Value = 4579.852
$ws.Range("K134").Value = 5495.6469
$ws.Range("L134").Value = 13739.556
$ws.Range("M134").Value = -2960.6469
$ws.Range("N134").Value = -18809.556

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 10500000
$ws.Range("J4").Value = 1000000
$ws.Range("L4").Value = 1000000
$ws.Range("N4").Value = -1000224
$ws.Range("H16").Value = 11415.1
$ws.Range("I16").Value = 14740
$ws.Range("J16").Value = 3657
$ws.Range("K16").Value = 14740
$ws.Range("L16").Value = 3657
$ws.Range("M16").Value = -14453
$ws.Range("N16").Value = -4231
$ws.Range("H31").Value = 4282.178
$ws.Range("I31").Value = 1395.2559
$ws.Range("J31").Value = 8420.1
$ws.Range("K31").Value = 1395.2559
$ws.Range("L31").Value = 8420.1
$ws.Range("M31").Value = -1100.2559
$ws.Range("N31").Value = -9010.1
$ws.Range("H34").Value = 4282.178
$ws.Range("I34").Value = 1395.2559
$ws.Range("J34").Value = 8420.1
$ws.Range("K34").Value = 1395.2559
$ws.Range("L34").Value = 8420.1
$ws.Range("M34").Value = -1193.2559
$ws.Range("N34").Value = -8824.1
$ws.Range("H36").Value = 20000
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 20000
$ws.Range("K36").Value = 0
$ws.Range("L36").ClearContents()
$ws.Range("M36").Value = 20000
$ws.Range("N36").Value = -20776
$ws.Range("H40").Value = 20000
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 20000
$ws.Range("K40").Value = 0
$ws.Range("L40").ClearContents()
$ws.Range("M40").Value = 20000
$ws.Range("N40").Value = -20320
$ws.Range("H99").Value = 2423.8125
$ws.Range("I99").Value = 1313.1428
$ws.Range("J99").Value = 3287.6667
$ws.Range("K99").Value = 1313.1428
$ws.Range("L99").Value = 3287.6667
$ws.Range("M99").Value = 184.8571999999999
$ws.Range("N99").Value = -6283.6667
$ws.Range("H105").Value = 1388.9615
$ws.Range("I105").Value = 1367.1052
$ws.Range("J105").Value = 1448.2858
$ws.Range("K105").Value = 1367.1052
$ws.Range("L105").Value = 1448.2858
$ws.Range("M105").Value = 379.8948
$ws.Range("N105").Value = -4942.2858
$ws.Range("H113").Value = 11415.1
$ws.Range("I113").Value = 14740
$ws.Range("J113").Value = 3657
$ws.Range("K113").Value = 14740
$ws.Range("L113").Value = 3657
$ws.Range("M113").Value = -12570
$ws.Range("N113").Value = -7997
$ws.Range("H126").Value = 2423.8125
$ws.Range("I126").Value = 1313.1428
$ws.Range("J126").Value = 3287.6667
$ws.Range("K126").Value = 3939.4284
$ws.Range("L126").Value = 9863.000100000001
$ws.Range("M126").Value = -1469.4284
$ws.Range("N126").Value = -14803.0001

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 5000147.5
$ws.Range("I4").Value = 5000147.5
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 15000442.5
$ws.Range("L4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -15000330.5
$ws.Range("H113").Value = 62502024
$ws.Range("I113").Value = 166668880
$ws.Range("J113").Value = 1910
$ws.Range("K113").Value = 500006640
$ws.Range("L113").Value = 5730
$ws.Range("M113").Value = -500004470
$ws.Range("N113").Value = -10070

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2886.9333
$ws.Range("I80").Value = 2670.5
$ws.Range("J80").Value = 3319.8
$ws.Range("K80").Value = 2670.5
$ws.Range("L80").Value = 3319.8
$ws.Range("M80").Value = -1672.5
$ws.Range("N80").Value = -5315.8
$ws.Range("H83").Value = 2886.9333
$ws.Range("I83").Value = 2670.5
$ws.Range("J83").Value = 3319.8
$ws.Range("K83").Value = 13352.5
$ws.Range("L83").Value = 16599
$ws.Range("M83").Value = -8360.5
$ws.Range("N83").Value = -26583
$ws.Range("H113").Value = 126538.875
$ws.Range("I113").Value = 200982.2
$ws.Range("J113").Value = 2466.6667
$ws.Range("K113").Value = 200982.2
$ws.Range("L113").Value = 2466.6667
$ws.Range("M113").Value = -198812.2
$ws.Range("N113").Value = -6806.6667
$ws.Range("H126").Value = 3167.375
$ws.Range("I126").Value = 1653.1111
$ws.Range("J126").Value = 5114.2856
$ws.Range("K126").Value = 4959.3333
$ws.Range("L126").Value = 15342.8568
$ws.Range("M126").Value = -2489.3333
$ws.Range("N126").Value = -20282.8568
$ws.Range("H132").Value = 6473.48
$ws.Range("I132").Value = 7135.1904
$ws.Range("K132").Value = 21405.5712
$ws.Range("M132").Value = -18875.5712

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 6250
$ws.Range("J2").Value = 10000
$ws.Range("L2").Value = 10000
$ws.Range("N2").Value = -10224
$ws.Range("H40").Value = 2960.9714
$ws.Range("I40").Value = 4824
$ws.Range("J40").Value = 1718.9524
$ws.Range("K40").Value = 4824
$ws.Range("L40").Value = 1718.9524
$ws.Range("M40").Value = -4688
$ws.Range("N40").Value = -1990.9524

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2049.875
$ws.Range("I81").Value = 499.75
$ws.Range("J81").Value = 3600
$ws.Range("K81").Value = 999.5
$ws.Range("L81").Value = 7200
$ws.Range("M81").Value = 61.5
$ws.Range("N81").Value = -9322
$ws.Range("H84").Value = 2049.875
$ws.Range("I84").Value = 499.75
$ws.Range("J84").Value = 3600
$ws.Range("K84").Value = 4997.5
$ws.Range("L84").Value = 36000
$ws.Range("M84").Value = 306.5
$ws.Range("N84").Value = -46608
$ws.Range("H107").Value = 1017.53845
$ws.Range("I107").Value = 1459.75
$ws.Range("J107").Value = 310
$ws.Range("K107").Value = 4379.25
$ws.Range("L107").Value = 930
$ws.Range("M107").Value = -2459.25
$ws.Range("N107").Value = -4770
